$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C6"   = -12.998
    "B7"   = 5.789
    "A9"   = -21.864
    "B12"  = 5.662000000000001
    "A13"  = -21.95
    "B14"  = 6.044
    "C15"  = -12.953
    "A16"  = -21.886
    "A18"  = -21.985
    "B19"  = 7.859999999999999
    "A20"  = -20.846
    "A26"  = -21.53299999999999
    "B26"  = 6.191
    "A27"  = -21.44
    "B27"  = 5.829
    "C28"  = -12.941
    "A29"  = -21.72
    "B29"  = 6.164
    "C33"  = -11.292
    "A35"  = -20.452
    "C35"  = -12.624
    "A36"  = -21.089
    "B37"  = 7.154999999999999
    "B38"  = 5.98
    "C38"  = -12.373
    "C43"  = -12.774
    "C44"  = -12.229
    "A45"  = -21.473
    "C45"  = -13.04
    "B47"  = 5.534000000000001
    "C47"  = -13.31
    "B51"  = 5.816
    "C51"  = -11.76
    "B52"  = 6.2
    "C54"  = -13.376
    "A55"  = -21.707
    "B55"  = 6.392
    "A57"  = -21.337
    "C57"  = -13.271
    "C62"  = -13.61
    "C63"  = -12.05
    "C67"  = -11.224
    "A69"  = -21.697
    "B69"  = 5.934
    "B70"  = 6.08
    "C70"  = -11.114
    "A76"  = -21.706
    "B76"  = 6.284999999999998
    "A78"  = -20.427
    "B81"  = 5.513
    "C81"  = -12.513
    "A82"  = -21.875
    "A83"  = -20.659
    "B83"  = 7.124
    "C88"  = -13.192
    "A93"  = -21.749
    "B94"  = 6.616999999999999
    "C96"  = -12.998
    "A97"  = -21.878
    "C99"  = -12.729
    "B100" = 6.127
    "B102" = 6.948
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
